# Auto-generated script to apply market-price data updates to Coeurl_Profits workbook
# Sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 714309.3
$ws.Range("I6").Value = 833353.8
$ws.Range("K6").Value = 2500061.4
$ws.Range("M6").Value = -2499949.4
$ws.Range("H19").Value = 1427.4667
$ws.Range("J19").Value = 1458
$ws.Range("L19").Value = 1458
$ws.Range("N19").Value = -1808
$ws.Range("H40").Value = 2463.75
$ws.Range("J40").Value = 2463.75
$ws.Range("L40").Value = 2463.75
$ws.Range("N40").Value = -2813.75
$ws.Range("H64").Value = 14357.143
$ws.Range("J64").Value = 15883.333
$ws.Range("L64").Value = 15883.333
$ws.Range("N64").Value = -16379.333
$ws.Range("H67").Value = 14357.143
$ws.Range("J67").Value = 15883.333
$ws.Range("L67").Value = 15883.333
$ws.Range("N67").Value = -17599.333
$ws.Range("H70").Value = 5100.75
$ws.Range("J70").Value = 8749.5
$ws.Range("L70").Value = 26248.5
$ws.Range("N70").Value = -26788.5
$ws.Range("H73").Value = 5100.75
$ws.Range("J73").Value = 8749.5
$ws.Range("L73").Value = 26248.5
$ws.Range("N73").Value = -28120.5
$ws.Range("H88").Value = 3354.5557
$ws.Range("I88").Value = 2999.5
$ws.Range("J88").Value = 3456
$ws.Range("K88").Value = 2999.5
$ws.Range("L88").Value = 3456
$ws.Range("M88").Value = -2593.5
$ws.Range("N88").Value = -4268
$ws.Range("H91").Value = 3354.5557
$ws.Range("I91").Value = 2999.5
$ws.Range("J91").Value = 3456
$ws.Range("K91").Value = 2999.5
$ws.Range("L91").Value = 3456
$ws.Range("M91").Value = -1595.5
$ws.Range("N91").Value = -6264
$ws.Range("H93").Value = 8974.75
$ws.Range("J93").Value = 8974.75
$ws.Range("L93").Value = 8974.75
$ws.Range("N93").Value = -13966.75
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
$ws.Range("H107").Value = 338.14285
$ws.Range("I107").Value = 338.14285
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 338.14285
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1581.85715
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 20066.666
$ws.Range("H32").Value = 4426.186
$ws.Range("I32").Value = 4157
$ws.Range("K32").Value = 4157
$ws.Range("M32").Value = -3870
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = ""
$ws.Range("H104").Value = 95120.5
$ws.Range("J104").Value = 95120.5
$ws.Range("L104").Value = 95120.5
$ws.Range("N104").Value = -102108.5
$ws.Range("H121").Value = 102999.5
$ws.Range("J121").Value = 102999.5
$ws.Range("L121").Value = 102999.5
$ws.Range("N121").Value = -106493.5
$ws.Range("H132").Value = 3095.3572
$ws.Range("I132").Value = 2768.3103
$ws.Range("K132").Value = 8304.930899999999
$ws.Range("M132").Value = -5774.930899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 554.5714
$ws.Range("I8").Value = 501.4
$ws.Range("J8").Value = 687.5
$ws.Range("K8").Value = 501.4
$ws.Range("L8").Value = 687.5
$ws.Range("M8").Value = -361.4
$ws.Range("N8").Value = -967.5
$ws.Range("H100").Value = 37166.668
$ws.Range("J100").Value = 37166.668
$ws.Range("L100").Value = 37166.668
$ws.Range("N100").Value = -39330.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 203.14285
$ws.Range("I5").Value = 149.45454
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 149.45454
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -37.45454000000001
$ws.Range("N5").Value = -624
$ws.Range("H20").Value = 147000
$ws.Range("J20").Value = 147000
$ws.Range("L20").Value = 147000
$ws.Range("N20").Value = -147472
$ws.Range("H22").Value = 643
$ws.Range("I22").Value = 640.8333
$ws.Range("J22").Value = 649.5
$ws.Range("K22").Value = 640.8333
$ws.Range("L22").Value = 649.5
$ws.Range("M22").Value = -290.8333
$ws.Range("N22").Value = -1349.5
$ws.Range("H25").Value = 1480
$ws.Range("I25").Value = 1100
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 1100
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -926
$ws.Range("N25").Value = -3348
$ws.Range("H30").Value = 147000
$ws.Range("J30").Value = 147000
$ws.Range("L30").Value = 147000
$ws.Range("N30").Value = -147182
$ws.Range("H86").Value = 5839.2
$ws.Range("I86").Value = 5122.5
$ws.Range("J86").Value = 6317
$ws.Range("K86").Value = 5122.5
$ws.Range("L86").Value = 6317
$ws.Range("M86").Value = -3999.5
$ws.Range("N86").Value = -8563
$ws.Range("H89").Value = 5839.2
$ws.Range("I89").Value = 5122.5
$ws.Range("J89").Value = 6317
$ws.Range("K89").Value = 25612.5
$ws.Range("L89").Value = 31585
$ws.Range("M89").Value = -19996.5
$ws.Range("N89").Value = -42817
$ws.Range("H110").Value = 250000
$ws.Range("J110").Value = 250000
$ws.Range("L110").Value = 250000
$ws.Range("N110").Value = -258180
$ws.Range("H122").Value = 1460.2142
$ws.Range("I122").Value = 1460.2142
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4380.642599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1930.642599999999
$ws.Range("N122").Value = ""
$ws.Range("H128").Value = 147000
$ws.Range("J128").Value = 147000
$ws.Range("L128").Value = 147000
$ws.Range("N128").Value = -156960
$ws.Range("H131").Value = 13001
$ws.Range("I131").Value = 13001
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 13001
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -7961
$ws.Range("N131").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 62783.75
$ws.Range("I9").Value = 117.5
$ws.Range("J9").Value = 125450
$ws.Range("K9").Value = 352.5
$ws.Range("L9").Value = 376350
$ws.Range("M9").Value = -128.5
$ws.Range("N9").Value = -376798
$ws.Range("H131").Value = 25558.674
$ws.Range("I131").Value = 334241.66
$ws.Range("J131").Value = 2407.45
$ws.Range("K131").Value = 1002724.98
$ws.Range("L131").Value = 7222.349999999999
$ws.Range("M131").Value = -997684.98
$ws.Range("N131").Value = -17302.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2339.5454
$ws.Range("I3").Value = 4124.5
$ws.Range("K3").Value = 4124.5
$ws.Range("M3").Value = -4008.5
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50346
$ws.Range("H104").Value = 32000
$ws.Range("J104").Value = 32000
$ws.Range("L104").Value = 32000
$ws.Range("N104").Value = -38988
$ws.Range("H114").Value = 26427
$ws.Range("J114").Value = 26427
$ws.Range("L114").Value = 26427
$ws.Range("N114").Value = -35105
$ws.Range("H132").Value = 3434
$ws.Range("I132").Value = 3467.875
$ws.Range("K132").Value = 10403.625
$ws.Range("M132").Value = -7873.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2213.6365
$ws.Range("I22").Value = 1916.8334
$ws.Range("K22").Value = 1916.8334
$ws.Range("M22").Value = -1621.8334
$ws.Range("H27").Value = 2213.6365
$ws.Range("I27").Value = 1916.8334
$ws.Range("K27").Value = 1916.8334
$ws.Range("M27").Value = -1809.8334
$ws.Range("H55").Value = 883.6923
$ws.Range("I55").Value = 689.8
$ws.Range("J55").Value = 1004.875
$ws.Range("K55").Value = 689.8
$ws.Range("L55").Value = 1004.875
$ws.Range("M55").Value = -516.8
$ws.Range("N55").Value = -1350.875
$ws.Range("H100").Value = 1593
$ws.Range("I100").Value = 1593
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1593
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1052
$ws.Range("N100").Value = ""
$ws.Range("H132").Value = 4893.4814
$ws.Range("I132").Value = 4436.75
$ws.Range("J132").Value = 6198.4287
$ws.Range("K132").Value = 13310.25
$ws.Range("L132").Value = 18595.2861
$ws.Range("M132").Value = -10780.25
$ws.Range("N132").Value = -23655.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 204285.58
$ws.Range("I2").Value = 250000
$ws.Range("K2").Value = 250000
$ws.Range("M2").Value = -249888
$ws.Range("H11").Value = 169993.17
$ws.Range("J11").Value = 3991.8
$ws.Range("L11").Value = 3991.8
$ws.Range("N11").Value = -4275.8
$ws.Range("H31").Value = 17
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""
$ws.Range("H81").Value = 11435.5
$ws.Range("J81").Value = 5272.727
$ws.Range("L81").Value = 10545.454
$ws.Range("N81").Value = -12667.454
$ws.Range("H84").Value = 11435.5
$ws.Range("J84").Value = 5272.727
$ws.Range("L84").Value = 52727.27
$ws.Range("N84").Value = -63335.27
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H122").Value = 2689.147
$ws.Range("I122").Value = 2558.5173
$ws.Range("K122").Value = 7675.5519
$ws.Range("M122").Value = -5225.5519
$ws.Range("H136").Value = 2345
$ws.Range("I136").Value = 1862.174
$ws.Range("K136").Value = 5586.522
$ws.Range("M136").Value = -3036.522
